$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(10, 1).Value = 9822.0499999999993
$ws.Cells.Item(10, 2).Value = 9895.27
$ws.Cells.Item(10, 3).Value = 286
$ws.Cells.Item(10, 4).Value = 283.87
$ws.Cells.Item(10, 5).Value = $false
$ws.Cells.Item(10, 6).Value = -0.74
$ws.Cells.Item(10, 7).Value = 42612.673032407409
$ws.Range("G10").NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(10, 8).Value = $false
